$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.998.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "'2.215.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'289.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.74%  "
$ws.Range("D6").Value = "'87.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.83%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +0.88%  "
$ws.Range("D10").Value = "'30.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.13%  "
$ws.Range("D11").Value = "'0.0778"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.82%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.110"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.60%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.46%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "'2.557.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("B15").Value = "Chainlink"
$ws.Range("C15").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D15").Value = "'13.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "'2.216.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.12%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "'0.729"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "'39.932.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").Value = "'11.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +11.92%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "'0.0₃0885"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'5.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").Value = "'65.56"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").Value = "'234.93"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "'2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("B26").Value = "ImmutableX"
$ws.Range("C26").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D26").Value = "'1.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.91%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'2.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.05%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'22.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").Value = "'9.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "'153.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.17%  "
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'31.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.15%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'4.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.01%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0719"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.03%  "
$ws.Range("B35").Value = "WEMIXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'2.39"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").Value = "'2.83"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +6.21%  "
$ws.Range("B37").Value = "Stellar"
$ws.Range("C37").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D37").Value = "'0.111"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").Value = "'0.0997"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.69%  "
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").Value = "'15.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.00%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "'1.71"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.67%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "'2.106.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.58%  "
$ws.Range("D42").Value = "'3.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.04%  "
$ws.Range("B43").Value = "ApeXProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D43").Value = "'2.19"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.86%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0268"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.59%  "
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'9.93"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.37%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'17.68"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.50%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'2.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.02%  "
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "'2.430.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "'1.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.95%  "
$ws.Range("D50").Value = "'69.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'88.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.35%  "
